$d = $word.ActiveDocument

# Widen the first column of the "Parameter / Lower CI / Estimate / Upper CI /
# Method" results table (table 8) from 2632 dxa to 2833 dxa (131.6pt -> 141.65pt).
$t = $d.Tables.Item(8)
$t.Columns(1).Width = 141.65

# "Difference in means" -> "Difference in medians" (row label, vMerge origin cell)
$d.Content.Find.Execute("Difference in means", $true, $false, $false, $false, $false, $true, 1, $false, "Difference in medians", 2)

# Bootstrapping row for Difference in medians: Lower CI / Estimate / Upper CI
$d.Content.Find.Execute("4.13", $true, $false, $false, $false, $false, $true, 1, $false, "3.68", 2)
$d.Content.Find.Execute("4.20", $true, $false, $false, $false, $false, $true, 1, $false, "3.73", 2)
$d.Content.Find.Execute("4.27", $true, $false, $false, $false, $false, $true, 1, $false, "3.79", 2)

# Bayesian Estimate row for Difference in medians: Lower CI / Estimate / Upper CI
$d.Content.Find.Execute("3.88", $true, $false, $false, $false, $false, $true, 1, $false, "3.81", 2)
$d.Content.Find.Execute("4.09", $true, $false, $false, $false, $false, $true, 1, $false, "3.92", 2)
$d.Content.Find.Execute("4.23", $true, $false, $false, $false, $false, $true, 1, $false, "4.02", 2)

# Percentage change, Bootstrapping row: Lower CI / Estimate / Upper CI
$d.Content.Find.Execute("30.77", $true, $false, $false, $false, $false, $true, 1, $false, "27.39", 2)
$d.Content.Find.Execute("31.28", $true, $false, $false, $false, $false, $true, 1, $false, "27.80", 2)
$d.Content.Find.Execute("31.80", $true, $false, $false, $false, $false, $true, 1, $false, "28.26", 2)

# Percentage change, Bayesian Estimate row: Lower CI / Estimate / Upper CI
$d.Content.Find.Execute("28.92", $true, $false, $false, $false, $false, $true, 1, $false, "28.38", 2)
$d.Content.Find.Execute("30.50", $true, $false, $false, $false, $false, $true, 1, $false, "29.23", 2)
$d.Content.Find.Execute("31.52", $true, $false, $false, $false, $false, $true, 1, $false, "29.98", 2)

# Effect size, Bayesian Estimate row: Lower CI / Estimate / Upper CI
$d.Content.Find.Execute("1.03", $true, $false, $false, $false, $false, $true, 1, $false, "1.41", 2)
$d.Content.Find.Execute("1.05", $true, $false, $false, $false, $false, $true, 1, $false, "1.43", 2)
$d.Content.Find.Execute("1.07", $true, $false, $false, $false, $false, $true, 1, $false, "1.46", 2)
